# Rename the two "audit" sheets to "total", and move the active/selected
# tab from by_course_KPI_audit (index 3) to the last sheet,
# by_course_demographic_audit / by_course_demographic_total (index 5).

$wb = $excel.ActiveWorkbook

$kpiSheet = $wb.Worksheets.Item("by_course_KPI_audit")
$kpiSheet.Name = "by_course_KPI_total"

$demoSheet = $wb.Worksheets.Item("by_course_demographic_audit")
$demoSheet.Name = "by_course_demographic_total"

# Make the renamed demographic sheet the active/selected tab, matching the
# workbook's new activeTab=5 / tabSelected placement.
$demoSheet.Activate()
